$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("L_ERSTT_2", "Q_EUROSTAT", "Erwerbstätigenquoten nach Geschlecht, Alter und Staatsangehörigkeit", "Employment rates by sex, age and citizenship", "https://ec.europa.eu/eurostat/databrowser/view/LFSA_ERGAN__custom_6067281/default/table?lang=de", "https://ec.europa.eu/eurostat/databrowser/view/LFSA_ERGAN__custom_6067281/default/table?lang=en"),
    @("L_ERSTT_3", "Q_EUROSTAT", "Gini-Koeffizient des verfügbaren Äquivalenzeinkommens vor Sozialleistungen - EU-SILC Erhebung - Eurostat-Tabelle  [ilc_di12c ]", "Gini coefficient of equivalised disposable income before social transfers - Eurostat table  [ilc_di12c ]", "https://ec.europa.eu/eurostat/databrowser/view/ILC_DI12C/default/table?lang=de&category=livcon.ilc.ilc_ie.ilc_iei", "https://ec.europa.eu/eurostat/databrowser/view/ILC_DI12C/default/table?lang=en&category=livcon.ilc.ilc_ie.ilc_iei"),
    @("L_ERSTT_4", "Q_EUROSTAT", "Rate der (erheblichen) materiellen Deprivation", "(Severe) Material deprivation rate", "https://ec.europa.eu/eurostat/databrowser/view/ILC_SIP8__custom_5385803/default/table?lang=de", "https://ec.europa.eu/eurostat/databrowser/view/ILC_SIP8__custom_5385803/default/table?lang=en"),
    @("L_ERSTT_5", "Q_EUROPEANCOMMISSION", "STECF - Berichte (nur auf Englisch verfügbar)", "STECF - Reports", "", "https://stecf.jrc.ec.europa.eu/reports/cfp-monitoring"),
    @("L_ERSTT_6", "Q_EUROSTAT", "Gini-Koeffizient des verfügbaren Äquivalenzeinkommens  - EU-SILC Erhebung - Eurostat-Tabelle  [ilc_di12 ]", "Gini coefficient of equivalised disposable income - Eurostat table  [ilc_di12c]", "https://ec.europa.eu/eurostat/databrowser/view/ILC_DI12/default/table?lang=de&category=livcon.ilc.ilc_ie.ilc_iei", "https://ec.europa.eu/eurostat/databrowser/view/ILC_DI12/default/table?lang=en&category=livcon.ilc.ilc_ie.ilc_iei"),
    @("L_ERSTT_7", "Q_EUROSTAT", "Rate der materiellen und sozialen Deprivation - Eurostat-Tabelle [ilc_mdsd07]", "Material and social deprivation rate - Eurostat table [ilc_mdsd07]", "https://ec.europa.eu/eurostat/databrowser/view/ILC_MDSD07/default/table?lang=de", "https://ec.europa.eu/eurostat/databrowser/view/ILC_MDSD07/default/table?lang=en"),
    @("L_ERSTT_8", "Q_EUROSTAT", "Rate der erheblichen materiellen Deprivation - Eurostat-Tabelle [ilc_mddd11]", "Severe material deprivation rate - Eurostat table [ilc_mddd11]", "https://ec.europa.eu/eurostat/databrowser/view/ILC_MDDD11/default/table?lang=de&category=livcon.ilc.ilc_md.ilc_mddd", "https://ec.europa.eu/eurostat/databrowser/view/ILC_MDDD11/default/table?lang=en&category=livcon.ilc.ilc_md.ilc_mddd"),
    @("L_ERSTT_9", "Q_EUROSTAT", "Bevölkerung nach Bildungsabschluss (nur Tertiärbereich) - Eurostat-Tabelle [edat_lfs_9912]", "Population by educational attainment level (tertiary education only) - Eurostat table [edat_lfs_9912]", "https://ec.europa.eu/eurostat/databrowser/view/EDAT_LFS_9912__custom_3201014/default/table?lang=de", "https://ec.europa.eu/eurostat/databrowser/view/EDAT_LFS_9912__custom_3201014/default/table?lang=en"),
    @("L_EZB_1", "Q_EZB", "Household Finance and Consumption Network (HFCN)", "Household Finance and Consumption Network (HFCN)", "https://www.ecb.europa.eu/pub/economic-research/research-networks/html/researcher_hfcn.en.html", "https://www.ecb.europa.eu/pub/economic-research/research-networks/html/researcher_hfcn.en.html"),
    @("L_FIDAR_1", "Q_FIDAR", "Studie zum WoB-Index", "Study on the WoB Index (only available in German)", "https://www.fidar.de/wob-indizes-studien/wob-index-185/studie-zum-wob-index-185.html", "https://www.fidar.de/wob-indizes-studien/wob-index-185/studie-zum-wob-index-185.html"),
    @("L_GBE_1", "Q_DESTATIS", "Vorzeitige Sterblichkeit", "Premature mortality", "https://www.gbe-bund.de/gbe/pkg_isgbe5.prc_menu_olap?p_uid=gast&p_aid=52889592&p_sprache=D&p_help=2&p_indnr=562&p_indsp=3194&p_ityp=H&p_fid=", "https://www.gbe-bund.de/gbe/pkg_isgbe5.prc_menu_olap?p_uid=gast&p_aid=80721089&p_sprache=E&p_help=2&p_indnr=562&p_version=1&p_ansnr=95662760"),
    @("L_GIZ_1", "Q_BMZ", "Mitglieder des Bündnises für nachhaltige Textilien", "Members of the Partnership for Sustainable Textiles", "https://www.textilbuendnis.com/uebersicht/", "https://www.textilbuendnis.com/en/uebersicht/"),
    @("L_IFEU_1", "Q_IFEU", "TREMOD", "TREMOD", "https://www.ifeu.de/methoden-tools/modelle/tremod/", "https://www.ifeu.de/en/methods-tools/models/tremod/"),
    @("L_LAK_1", "Q_LAKEB", "Energieindikatoren allgemein", "Energy indicators general (only available in German)", "https://www.lak-energiebilanzen.de/eingabe-dynamisch/?a=i100", "https://www.lak-energiebilanzen.de/eingabe-dynamisch/?a=i100"),
    @("L_LAK_2", "Q_LAKEB", "Energieindikatoren: Erneuerbare Energieträger", "Energy indicators: Renewable energy sources (only available in German)", "https://www.lak-energiebilanzen.de/eingabe-dynamisch/?a=i290", "https://www.lak-energiebilanzen.de/eingabe-dynamisch/?a=i290"),
    @("L_LAK_3", "Q_LAKEB", "Anteil Strom aus erneuerbaren Energiequellen am Bruttostromverbrauch", "Share of electricity from renewable energy sources in gross electricity consumption (only available in German)", "https://www.lak-energiebilanzen.de/eingabe-dynamisch/?a=i200", "https://www.lak-energiebilanzen.de/eingabe-dynamisch/?a=i200"),
    @("L_LANUV_2", "Q_LIKI", "Flächenverbrauch", "Land consumption (only available in German)", "https://www.liki.nrw.de/ressourcen-und-effizienz/d1-flaechenverbrauch", ""),
    @("L_LIKI_1", "Q_LIKI", "Nitrat im Grundwasser", "Nitrate in groundwater (only available in German)", "https://www.liki.nrw.de/umwelt-und-gesundheit/c5-nitrat-im-grundwasser", ""),
    @("L_OECD_1", "Q_OECD", "OECD - Gesamtströme nach Geber (nicht auf Deutsch verfügbar)", "OECD - Total flows by donor", "", "https://stats.oecd.org/Index.aspx?DataSetCode=TABLE1"),
    @("L_RDB_1", "Q_DESTATIS", "Bevölkerung", "Population (only available in German)", "https://www.regionalstatistik.de/genesis/online?operation=previous&levelindex=0&step=0&titel=Tabellenaufbau&levelid=1668672879939&acceptscookies=false#abreadcrumb", ""),
    @("L_RDB_2", "Q_DESTATIS", "Siedlungsfläche", "Settlement area (only available in German)", "https://www.regionalstatistik.de/genesis//online?operation=table&code=33111-02-01-4&bypass=true&levelindex=1&levelid=1668501420853#abreadcrumb", ""),
    @("L_RDB_3", "Q_DESTATIS", "Verkehrsfläche", "Transport area (only available in German)", "https://www.regionalstatistik.de/genesis//online?operation=table&code=33111-03-01-4&bypass=true&levelindex=1&levelid=1668501420853#abreadcrumb", ""),
    @("L_SP_1", "Q_VWGDL", "Investitionen, Anlagevermögen (eigene Berechnung auf Basis der Volkswirtschaftlichen Gesamtrechnungen der Länder)", "Investments, fixed assets (own calculation based on the Volkswirtschaftliche Gesamtrechnungen der Länder - only available in German)", "https://www.statistikportal.de/de/vgrdl/ergebnisse-laenderebene/investitionen-anlagevermoegen", "https://www.statistikportal.de/de/vgrdl/ergebnisse-laenderebene/investitionen-anlagevermoegen"),
    @("L_SP_2", "Q_VWGDL", "Bruttoinlandsprodukt, Bruttowertschöpfung (eigene Berechnung auf Basis der Volkswirtschaftlichen Gesamtrechnungen)", "Gross domestic product, gross value added (own calculation based on the Volkswirtschaftliche Gesamtrechnungen der Länder  - only available in German)", "https://www.statistikportal.de/de/vgrdl/ergebnisse-laenderebene/bruttoinlandsprodukt-bruttowertschoepfung", "https://www.statistikportal.de/de/vgrdl/ergebnisse-laenderebene/bruttoinlandsprodukt-bruttowertschoepfung"),
    @("L_SP_3", "Q_VWGDL", "Treibhausgasemissionen ", "Greenhouse gas emissions (only available in German)", "https://www.statistikportal.de/de/ugrdl/ergebnisse/gase#alle-ergebnisse", "https://www.statistikportal.de/de/ugrdl/ergebnisse/gase#alle-ergebnisse"),
    @("L_SP_4", "Q_VWGDL", "Erwerbstätigenquoten", "Employment rates (only available in German)", "https://www.statistikportal.de/de/nachhaltigkeit/ergebnisse/ziel-8-menschenwuerdige-arbeit-und-wirtschaftswachstum", ""),
    @("L_SVWS", "Q_SVWS", "Forschung und Entwicklung in der Wirtschaft", "Research and development in the economy (only available in German)", "https://www.stifterverband.org/fue-facts-2020", ""),
    @("L_TI_1", "Q_TA", "Corruption Perceptions Index (nicht auf Deutsch verfügbar)", "Corruption Perceptions Index", "https://www.transparency.org/en/cpi/2020/index/nzl", "https://www.transparency.org/en/cpi/2020/index/nzl"),
    @("L_UBA_1", "Q_UBA", "Index der Luftschadstoff-Emissionen", "Index of air pollutant emissions", "https://www.umweltbundesamt.de/bild/index-der-luftschadstoff-emissionen", "https://www.umweltbundesamt.de/en/image/index-of-air-pollutant-emissions"),
    @("L_UBA_2", "Q_UBALAWA", "Messstellen an Flüssen mit Überschreitung des Orientierungswertes für Gesamtphosphor", "Sampling sites which exceeded the requirement for good status for total phosphorus in rivers", "https://www.umweltbundesamt.de/daten/umweltindikatoren/indikator-eutrophierung-von-fluessen-durch-phosphor", "https://www.umweltbundesamt.de/en/data/environmental-indicators/indicator-river-eutrophication-phosphorus"),
    @("L_UBA_3", "Q_UBA", "Nitrat im Grundwasser", "Nitrate in groundwater", "https://www.umweltbundesamt.de/daten/umweltindikatoren/indikator-nitrat-im-grundwasser", "https://www.umweltbundesamt.de/en/data/environmental-indicators/indicator-nitrate-in-groundwater"),
    @("L_UBA_4", "Q_UBA", "Marktanteile von Produkten mit staatlichen Umweltzeichen, nach Umsätzen gewichtet", "Weighted market shares by sales of products with official eco-labels", "https://www.umweltbundesamt.de/daten/umweltindikatoren/indikator-umweltfreundlicher-konsum", "https://www.umweltbundesamt.de/en/indicator-environmentally-friendly-consumption"),
    @("L_UBA_5", "Q_UBA", "Emission der von der UN-Klimarahmenkonvention abgedeckten Treibhausgase", "Emission of greenhouse gases covered by the UN Framework Convention on Climate", "https://www.umweltbundesamt.de/daten/umweltindikatoren/indikator-emission-von-treibhausgasen", "https://www.umweltbundesamt.de/en/data/environmental-indicators/indicator-greenhouse-gas-emissions"),
    @("L_UBA_6", "Q_UBA_1", "Eutrophierung von Nord- und Ostsee durch Stickstoff", "Eutrophication of the North Sea / Baltic Sea by nitrogen", "https://www.umweltbundesamt.de/daten/umweltindikatoren/indikator-eutrophierung-der-meere", "https://www.umweltbundesamt.de/en/data/environmental-indicators/indicator-eutrophication-of-the-north-sea-baltic-sea"),
    @("L_UBA_7", "Q_UBA", "Anteil der Fläche empfindlicher Land-Ökosysteme mit Überschreitung der Belastungsgrenzen für Eutrophierung", "Proportion of vulnerable ecosystems where critical loads for eutrophication are exceeded", "https://www.umweltbundesamt.de/daten/umweltindikatoren/indikator-eutrophierung-durch-stickstoff", "https://www.umweltbundesamt.de/en/data/environmental-indicators/indicator-nitrogen-eutrophication")
)

$startRow = 45
$i = 0
foreach ($item in $data) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
    $ws.Cells.Item($r, 6).Value = $item[5]
    $i = $i + 1
}

Write-Host "Updated $i rows"